$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in a new test-run column (K) for already existing test rows (3-45) ---
# Row 2 already has K2 = "OK" recorded from a previous run.
$ws.Range("K3:K45").Value = "OK"

# --- New bug reports discovered while fixing the "looping on buggy code" issue ---
$nl = [char]10

# Row 46: infinite loop / NameError bug - now fixed => OK
$text46 = "NameError: name 'aaa' is not defined зацикливается" + $nl + `
  "def func0(arg0):" + $nl + `
  "  if compare(1, ""<"", arg0):" + $nl + `
  "    pour(1, 2)" + $nl + `
  "func0(aaa)" + $nl

$ws.Range("B46").WrapText = $true
$ws.Range("B46").Value = $text46
$ws.Range("K46").Value = "OK"
$ws.Rows(46).RowHeight = 15

# Row 47: new idea/issue, still unresolved => "?"
$text47 = "подставлять значения аргументов в каунтеры и селекты? Опасно!!!"
$ws.Range("B47").WrapText = $true
$ws.Range("B47").Value = $text47
$ws.Range("K47").Value = "?"

# Row 48: new idea/issue, still unresolved => "?"
$text48 = "добавить возможность добавлять в роли аргументов строки, которыеи могут быть именаим переменных???"
$ws.Range("B48").WrapText = $true
$ws.Range("B48").Value = $text48
$ws.Range("K48").Value = "?"
$ws.Rows(48).RowHeight = 14.25

# --- Highlight the still-open ("?") items the same way the rest of the sheet does ---
foreach ($addr in @("K47", "K48")) {
    $rng = $ws.Range($addr)
    $fc = $rng.FormatConditions.Add(1, 3, '"?"')
    $fc.Font.Color = 26012
    $fc.Interior.Color = 10284031
}

# --- Leave the cursor where the author left it after adding the new rows ---
$ws.Range("B51").Select()
